$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new election candidate row (row 9) with roll no, name, dept, prog, position
$ws.Range("A9").Value = "B180010"
$ws.Range("B9").Value = "appukutan"
$ws.Range("C9").Value = "FE"
$ws.Range("D9").Value = "PHD"
$ws.Range("E9").Value = "Entho secretary"

# Leave the last edited cell selected, matching the recorded selection after data entry
$ws.Range("E9").Select()
